$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 16 runtimes: Part 1 and Part 2 timings (seconds)
$ws.Range("B20").Value = 0.0081190999771934003
$ws.Range("C20").Value = 1.6652642000117299

# Move the active selection to A20, as was left after entering the data
$ws.Range("A20").Select()
